$p = $ppt.ActivePresentation
$p.Slides.Item(27).Delete()
